$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Chad"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.4370383333333334
$ws.Cells.Item(2, 8).Value = 1.311115
$ws.Cells.Item(2, 9).Value = 0.04569839301109439
$ws.Cells.Item(2, 10).Value = 0.04569839301109439
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 153.5290173333333
$ws.Cells.Item(2, 14).Value = 460.587052
$ws.Cells.Item(2, 15).Value = 0.3172206968818489
$ws.Cells.Item(2, 16).Value = 0.317220696881849
$ws.Cells.Item(2, 17).Value = 67.09806585366444
$ws.Cells.Item(2, 18).Value = 603.88259268298
$ws.Cells.Item(2, 19).Value = 0.01449647607735998
$ws.Cells.Item(2, 20).Value = 0.01449647607735998

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Chad"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.4370383333333334
$ws.Cells.Item(3, 8).Value = 1.311115
$ws.Cells.Item(3, 9).Value = 0.04569839301109439
$ws.Cells.Item(3, 10).Value = 0.04569839301109439
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 168.7997026666667
$ws.Cells.Item(3, 14).Value = 506.3991080000001
$ws.Cells.Item(3, 15).Value = 0.3487728915577651
$ws.Cells.Item(3, 16).Value = 0.3487728915577651
$ws.Cells.Item(3, 17).Value = 73.77194072060223
$ws.Cells.Item(3, 18).Value = 663.9474664854201
$ws.Cells.Item(3, 19).Value = 0.01593836067002255
$ws.Cells.Item(3, 20).Value = 0.01593836067002255

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Chad"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.4370383333333334
$ws.Cells.Item(4, 8).Value = 1.311115
$ws.Cells.Item(4, 9).Value = 0.04569839301109439
$ws.Cells.Item(4, 10).Value = 0.04569839301109439
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 68.09032333333333
$ws.Cells.Item(4, 14).Value = 204.27097
$ws.Cells.Item(4, 15).Value = 0.1406878008722904
$ws.Cells.Item(4, 16).Value = 0.1406878008722904
$ws.Cells.Item(4, 17).Value = 29.75808142572778
$ws.Cells.Item(4, 18).Value = 267.82273283155
$ws.Cells.Item(4, 19).Value = 0.006429206416128514
$ws.Cells.Item(4, 20).Value = 0.006429206416128516

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Chad"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.4370383333333334
$ws.Cells.Item(5, 8).Value = 1.311115
$ws.Cells.Item(5, 9).Value = 0.04569839301109439
$ws.Cells.Item(5, 10).Value = 0.04569839301109439
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 93.562673
$ws.Cells.Item(5, 14).Value = 280.688019
$ws.Cells.Item(5, 15).Value = 0.1933186106880956
$ws.Cells.Item(5, 16).Value = 0.1933186106880956
$ws.Cells.Item(5, 17).Value = 40.89047467013167
$ws.Cells.Item(5, 18).Value = 368.014272031185
$ws.Cells.Item(5, 19).Value = 0.008834349847583347
$ws.Cells.Item(5, 20).Value = 0.008834349847583347

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Chad"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.326553
$ws.Cells.Item(6, 8).Value = 18.979659
$ws.Cells.Item(6, 9).Value = 0.6615284823974669
$ws.Cells.Item(6, 10).Value = 0.6615284823974669
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 153.5290173333333
$ws.Cells.Item(6, 14).Value = 460.587052
$ws.Cells.Item(6, 15).Value = 0.3172206968818489
$ws.Cells.Item(6, 16).Value = 0.317220696881849
$ws.Cells.Item(6, 17).Value = 971.3094651972518
$ws.Cells.Item(6, 18).Value = 8741.785186775267
$ws.Cells.Item(6, 19).Value = 0.2098505261933163
$ws.Cells.Item(6, 20).Value = 0.2098505261933164

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Chad"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.326553
$ws.Cells.Item(7, 8).Value = 18.979659
$ws.Cells.Item(7, 9).Value = 0.6615284823974669
$ws.Cells.Item(7, 10).Value = 0.6615284823974669
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 168.7997026666667
$ws.Cells.Item(7, 14).Value = 506.3991080000001
$ws.Cells.Item(7, 15).Value = 0.3487728915577651
$ws.Cells.Item(7, 16).Value = 0.3487728915577651
$ws.Cells.Item(7, 17).Value = 1067.920265304908
$ws.Cells.Item(7, 18).Value = 9611.282387744172
$ws.Cells.Item(7, 19).Value = 0.2307232016535846
$ws.Cells.Item(7, 20).Value = 0.2307232016535846

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Chad"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 6.326553
$ws.Cells.Item(8, 8).Value = 18.979659
$ws.Cells.Item(8, 9).Value = 0.6615284823974669
$ws.Cells.Item(8, 10).Value = 0.6615284823974669
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 68.09032333333333
$ws.Cells.Item(8, 14).Value = 204.27097
$ws.Cells.Item(8, 15).Value = 0.1406878008722904
$ws.Cells.Item(8, 16).Value = 0.1406878008722904
$ws.Cells.Item(8, 17).Value = 430.77703935547
$ws.Cells.Item(8, 18).Value = 3876.99335419923
$ws.Cells.Item(8, 19).Value = 0.09306898740288327
$ws.Cells.Item(8, 20).Value = 0.09306898740288329

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Chad"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 6.326553
$ws.Cells.Item(9, 8).Value = 18.979659
$ws.Cells.Item(9, 9).Value = 0.6615284823974669
$ws.Cells.Item(9, 10).Value = 0.6615284823974669
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 93.562673
$ws.Cells.Item(9, 14).Value = 280.688019
$ws.Cells.Item(9, 15).Value = 0.1933186106880956
$ws.Cells.Item(9, 16).Value = 0.1933186106880956
$ws.Cells.Item(9, 17).Value = 591.929209556169
$ws.Cells.Item(9, 18).Value = 5327.362886005521
$ws.Cells.Item(9, 19).Value = 0.1278857671476826
$ws.Cells.Item(9, 20).Value = 0.1278857671476826

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Chad"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8278226666666667
$ws.Cells.Item(10, 8).Value = 2.483468
$ws.Cells.Item(10, 9).Value = 0.08656029157966813
$ws.Cells.Item(10, 10).Value = 0.08656029157966813
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 153.5290173333333
$ws.Cells.Item(10, 14).Value = 460.587052
$ws.Cells.Item(10, 15).Value = 0.3172206968818489
$ws.Cells.Item(10, 16).Value = 0.317220696881849
$ws.Cells.Item(10, 17).Value = 127.0948005395929
$ws.Cells.Item(10, 18).Value = 1143.853204856336
$ws.Cells.Item(10, 19).Value = 0.02745871601719836
$ws.Cells.Item(10, 20).Value = 0.02745871601719837

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Chad"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.8278226666666667
$ws.Cells.Item(11, 8).Value = 2.483468
$ws.Cells.Item(11, 9).Value = 0.08656029157966813
$ws.Cells.Item(11, 10).Value = 0.08656029157966813
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 168.7997026666667
$ws.Cells.Item(11, 14).Value = 506.3991080000001
$ws.Cells.Item(11, 15).Value = 0.3487728915577651
$ws.Cells.Item(11, 16).Value = 0.3487728915577651
$ws.Cells.Item(11, 17).Value = 139.7362199940605
$ws.Cells.Item(11, 18).Value = 1257.625979946544
$ws.Cells.Item(11, 19).Value = 0.03018988318832412
$ws.Cells.Item(11, 20).Value = 0.03018988318832412

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Chad"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.8278226666666667
$ws.Cells.Item(12, 8).Value = 2.483468
$ws.Cells.Item(12, 9).Value = 0.08656029157966813
$ws.Cells.Item(12, 10).Value = 0.08656029157966813
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 68.09032333333333
$ws.Cells.Item(12, 14).Value = 204.27097
$ws.Cells.Item(12, 15).Value = 0.1406878008722904
$ws.Cells.Item(12, 16).Value = 0.1406878008722904
$ws.Cells.Item(12, 17).Value = 56.36671303599556
$ws.Cells.Item(12, 18).Value = 507.3004173239601
$ws.Cells.Item(12, 19).Value = 0.01217797706520774
$ws.Cells.Item(12, 20).Value = 0.01217797706520775

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Chad"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.8278226666666667
$ws.Cells.Item(13, 8).Value = 2.483468
$ws.Cells.Item(13, 9).Value = 0.08656029157966813
$ws.Cells.Item(13, 10).Value = 0.08656029157966813
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 93.562673
$ws.Cells.Item(13, 14).Value = 280.688019
$ws.Cells.Item(13, 15).Value = 0.1933186106880956
$ws.Cells.Item(13, 16).Value = 0.1933186106880956
$ws.Cells.Item(13, 17).Value = 77.45330146332134
$ws.Cells.Item(13, 18).Value = 697.0797131698921
$ws.Cells.Item(13, 19).Value = 0.01673371530893791
$ws.Cells.Item(13, 20).Value = 0.01673371530893791

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Chad"
$ws.Cells.Item(14, 3).Value = "Itgb1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 1.972124333333333
$ws.Cells.Item(14, 8).Value = 5.916373
$ws.Cells.Item(14, 9).Value = 0.2062128330117706
$ws.Cells.Item(14, 10).Value = 0.2062128330117706
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 153.5290173333333
$ws.Cells.Item(14, 14).Value = 460.587052
$ws.Cells.Item(14, 15).Value = 0.3172206968818489
$ws.Cells.Item(14, 16).Value = 0.317220696881849
$ws.Cells.Item(14, 17).Value = 302.7783109558217
$ws.Cells.Item(14, 18).Value = 2725.004798602396
$ws.Cells.Item(14, 19).Value = 0.0654149785939742
$ws.Cells.Item(14, 20).Value = 0.0654149785939742

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Chad"
$ws.Cells.Item(15, 3).Value = "Itgb1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 1.972124333333333
$ws.Cells.Item(15, 8).Value = 5.916373
$ws.Cells.Item(15, 9).Value = 0.2062128330117706
$ws.Cells.Item(15, 10).Value = 0.2062128330117706
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 168.7997026666667
$ws.Cells.Item(15, 14).Value = 506.3991080000001
$ws.Cells.Item(15, 15).Value = 0.3487728915577651
$ws.Cells.Item(15, 16).Value = 0.3487728915577651
$ws.Cells.Item(15, 17).Value = 332.8940010883649
$ws.Cells.Item(15, 18).Value = 2996.046009795285
$ws.Cells.Item(15, 19).Value = 0.07192144604583378
$ws.Cells.Item(15, 20).Value = 0.07192144604583378

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Chad"
$ws.Cells.Item(16, 3).Value = "Itgb1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 1.972124333333333
$ws.Cells.Item(16, 8).Value = 5.916373
$ws.Cells.Item(16, 9).Value = 0.2062128330117706
$ws.Cells.Item(16, 10).Value = 0.2062128330117706
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 68.09032333333333
$ws.Cells.Item(16, 14).Value = 204.27097
$ws.Cells.Item(16, 15).Value = 0.1406878008722904
$ws.Cells.Item(16, 16).Value = 0.1406878008722904
$ws.Cells.Item(16, 17).Value = 134.2825835102011
$ws.Cells.Item(16, 18).Value = 1208.54325159181
$ws.Cells.Item(16, 19).Value = 0.02901162998807084
$ws.Cells.Item(16, 20).Value = 0.02901162998807085

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Chad"
$ws.Cells.Item(17, 3).Value = "Itgb1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 1.972124333333333
$ws.Cells.Item(17, 8).Value = 5.916373
$ws.Cells.Item(17, 9).Value = 0.2062128330117706
$ws.Cells.Item(17, 10).Value = 0.2062128330117706
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 93.562673
$ws.Cells.Item(17, 14).Value = 280.688019
$ws.Cells.Item(17, 15).Value = 0.1933186106880956
$ws.Cells.Item(17, 16).Value = 0.1933186106880956
$ws.Cells.Item(17, 17).Value = 184.5172241150097
$ws.Cells.Item(17, 18).Value = 1660.655017035087
$ws.Cells.Item(17, 19).Value = 0.03986477838389175
$ws.Cells.Item(17, 20).Value = 0.03986477838389175
